# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly scraped totals, per commit "Update gh-pages to output
# generated at 8a634ce".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1230
$ws1.Range("F3").Value = 1989
$ws1.Range("F5").Value = 168
$ws1.Range("F6").Value = 405
$ws1.Range("F8").Value = 495
$ws1.Range("F9").Value = 121
$ws1.Range("F13").Value = 48
$ws1.Range("F15").Value = 3795
$ws1.Range("F17").Value = 825
$ws1.Range("F19").Value = 333
$ws1.Range("F20").Value = 708
$ws1.Range("F21").Value = 1248
$ws1.Range("F22").Value = 44
$ws1.Range("F25").Value = 77

# Sheet "全部类型" (sheet4): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1230
$ws4.Range("F7").Value = 1989
$ws4.Range("F9").Value = 168
$ws4.Range("F10").Value = 405
$ws4.Range("F12").Value = 495
$ws4.Range("F13").Value = 121
$ws4.Range("F17").Value = 48
$ws4.Range("F21").Value = 3795
$ws4.Range("F23").Value = 825
$ws4.Range("F25").Value = 333
$ws4.Range("F26").Value = 708
$ws4.Range("F27").Value = 1248
$ws4.Range("F28").Value = 44
$ws4.Range("F31").Value = 77
